$wb = $excel.ActiveWorkbook

# Sheet "Metadata": update URL, Version, Date values
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/ValueSet/MeltingPotVS"
$wsMeta.Range("B3").Value = "2.0.0"
$wsMeta.Range("B8").Value = "2026-01-15T15:23:39+00:00"

# Sheet "Include #1": update competence-code-system URI
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/competence-code-system"

# Sheet "Include #2": update type-carte-code-system URI
$wsInc2 = $wb.Worksheets.Item("Include #2")
$wsInc2.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/tde/CodeSystem/type-carte-code-system"
